$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test case TC-003 data (rows 7 & 8) ---
# Fill column A-D first for both rows, then E-H for both rows, to
# reproduce the original shared-string insertion order.
$ws.Range("A7").Value = "TC-003"
$ws.Range("B7").Value = "expTitle"
$ws.Range("C7").Value = "username"
$ws.Range("D7").Value = "password"

$ws.Range("A8").Value = "TC-003"
$ws.Range("B8").Value = "expTitle"
$ws.Range("C8").Value = "kiran2403"
$ws.Range("D8").Value = "kiran123"

$ws.Range("E7").Value = "confirm Password"
$ws.Range("F7").Value = "Full  Name"
$ws.Range("G7").Value = "email id"
$ws.Range("H7").Value = "captcha"

$ws.Range("E8").Value = "kiran123"
$ws.Range("F8").Value = "Kiran Kumar"
$ws.Range("G8").Value = "abcd@gmail.com"
$ws.Range("H8").Value = "Hyderabad"

# --- Hyperlink on the email address cell ---
$ws.Hyperlinks.Add($ws.Range("G8"), "mailto:abcd@gmail.com")

# --- Column widths (approximate AutoFit results for the new/changed columns) ---
$ws.Columns.Item(4).ColumnWidth = 49.333333333333336
$ws.Columns.Item(5).ColumnWidth = 25.166666666666668
$ws.Columns.Item(6).ColumnWidth = 15.0
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 15.333333333333334

# --- Final selection, as left by the author ---
$null = $ws.Range("I6").Select()
